$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packet")

$ws.Range("E2").Value = "Packet"
$ws.Range("F2").Value = "pk"
$ws.Range("E3").Value = "Packet"
$ws.Range("E4").Value = "Packet"
$ws.Range("E5").Value = "Packet"
$ws.Range("E6").Value = "Packet"
$ws.Range("E7").Value = "Packet"
$ws.Range("E8").Value = "Packet"
$ws.Range("E9").Value = "Packet"
$ws.Range("E10").Value = "Packet"
$ws.Range("E11").Value = "Packet"
$ws.Range("E12").Value = "Packet"

$ws.Range("F2").Select()
